# Update the "想去人数" (interested-count) values in column F across all
# four worksheets, per the source data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> @{ row = newValue }
$updates = @{
    "展览" = @{
        2  = 502
        3  = 1558
        4  = 811
        5  = 217
        6  = 57
        7  = 1118
        8  = 710
        9  = 765
        10 = 1381
        12 = 1018
        13 = 29
        14 = 61
        15 = 185
        17 = 432
        18 = 7
        20 = 290
        21 = 536
        22 = 556
        23 = 744
        24 = 231
        26 = 366
    }
    "演出" = @{
        3  = 990
        5  = 249
        7  = 138
        9  = 583
        10 = 74
    }
    "本地生活" = @{
        2 = 206
    }
    "全部类型" = @{
        2  = 502
        3  = 206
        4  = 1558
        6  = 811
        7  = 217
        8  = 990
        9  = 57
        10 = 1118
        11 = 710
        12 = 765
        13 = 1381
        15 = 1018
        16 = 29
        17 = 61
        18 = 185
        20 = 432
        21 = 7
        23 = 249
        25 = 290
        27 = 138
        28 = 138
        29 = 536
        30 = 556
        31 = 744
        32 = 231
        35 = 583
        36 = 74
        37 = 74
        39 = 366
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $updates[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $ws.Range("F$row").Value = $rowsForSheet[$row]
    }
}
